$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target stored width = 72; runtime adds ~0.8333 offset
# when converting ColumnWidth -> stored OOXML width, so compensate here)
$ws.Columns.Item(1).ColumnWidth = 71.16666666666667

# Update header row values
$ws.Range("A1").Value = "Page URL"
$ws.Range("B1").Value = "Test Case"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Comments"

# Update the URL in A2
$ws.Range("A2").Value = "https://www.alojamiento.io/property/mall-of-i-stanbul-3/BC-6975002/"
